$d = $word.ActiveDocument

# "Lower row of " and "J3" were two separate runs; merge them into one run
# with the combined text (no other formatting change).
$d.Content.Find.Execute("Lower row of J3", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Lower row of J3", 2)

# The footer "CLASSICV2 Additional Information v1 (Nov 2021)" text was split
# across four runs; collapse them into a single run carrying the new title.
foreach ($sec in $d.Sections) {
    $ftr = $sec.Footers.Item(1)
    if ($ftr.Exists) {
        $ftr.Range.Find.Execute("CLASSICV2 Additional Information v1 (Nov 2021)", `
                                 $true, $false, $false, $false, $false, $true, 1, $false, `
                                 "rosco_m68k Keyboard (r1) Errata (Jan 2024)", 2)
    }
}
